$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.249.11"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.856.67"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.45"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4745"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2747"
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "1.873.88"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07426"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.07"
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.983"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.33"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6329"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("D16").Value = "30.226.44"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.76"
$ws.Range("E18").Value = "  -4.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007306"
$ws.Range("E19").Value = "  -3.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "225.13"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "2.089.30"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.099"
$ws.Range("E23").Value = "  -3.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.024"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.41"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.217"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.79"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.856"
$ws.Range("E28").Value = "  -5.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1023"
$ws.Range("E29").Value = "  +9.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.382"
$ws.Range("E30").Value = "  -5.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.225"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.902"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04879"
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7263"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9999"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.686"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01916"
$ws.Range("E38").Value = "  +5.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.627"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8999"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.973"
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.69"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9941"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4096"
$ws.Range("E44").Value = "  -4.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.526"
$ws.Range("E45").Value = "  -6.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.039"
$ws.Range("E46").Value = "  -5.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.34"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.786"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.401"
$ws.Range("E50").Value = "  -5.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05588"
$ws.Range("E51").Value = "  -0.79%  "
